$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column AG (33rd column) -- this is the new
# DC_IDENTIFIER_TYPE column. Excel shifts every column from AG onward one
# position to the right (AG->AH, AH->AI, ... AY->AZ) and the new blank
# column inherits formatting from the old column it displaces.
$ws.Columns("AG:AG").Insert()

# Header for the newly inserted column.
$ws.Range("AG1").Value = "DC_IDENTIFIER_TYPE"

# Row 6 (ruimtereis02 relation block): add an ISSN identifier + its type.
$ws.Range("AG6").Value = "ISSN"
$ws.Range("AF6").Value = "0925-6229"

# Row 3 (ruimtereis01 relation block): add an identifier value + its type.
$ws.Range("AG3").Value = "ARCHIS-ZAAK-IDENTIFICATIE"
$ws.Range("AF3").Value = 6663

# The hidden _FilterDatabase defined name covers the header row and must
# grow by one column to keep including the new last column (AY, was AX).
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "=Sheet1!`$A`$1:`$AY`$9"

# Restore the selection to its documented post-edit state.
$ws.Range("AG3").Select()
